$wb = $excel.ActiveWorkbook

# Update the "Last Updated" timestamp on the Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 03:05 PM"

# Update the "1 Year" column (F) values on the Industry Analysis sheet
$ws = $wb.Worksheets.Item("Industry Analysis")
$ws.Range("F2").Value = 21.3
$ws.Range("F3").Value = -4.3927
$ws.Range("F4").Value = 35.9445
$ws.Range("F5").Value = -51.0482
$ws.Range("F6").Value = 57.2275
$ws.Range("F7").Value = -9.640700000000001
$ws.Range("F8").Value = -6.1449
$ws.Range("F9").Value = 36.9733
$ws.Range("F10").Value = -4.7026
$ws.Range("F11").Value = 46.5317
$ws.Range("F12").Value = -2.102
$ws.Range("F13").Value = 17.4681
$ws.Range("F14").Value = -33.0245
$ws.Range("F15").Value = 1.0205
$ws.Range("F16").Value = 2.0426
$ws.Range("F17").Value = -16.2411
$ws.Range("F18").Value = 7.4627
$ws.Range("F19").Value = -25.798
$ws.Range("F20").Value = 47.7485
$ws.Range("F21").Value = 19.5587
$ws.Range("F22").Value = 76.5603
$ws.Range("F23").Value = -54.2675
$ws.Range("F24").Value = -0.8811
$ws.Range("F25").Value = 4.8518
$ws.Range("F26").Value = 3.6831
$ws.Range("F27").Value = -34.0874
$ws.Range("F28").Value = -11.9893
$ws.Range("F29").Value = -12.994
$ws.Range("F30").Value = 25.5415
$ws.Range("F31").Value = 56.5088
$ws.Range("F32").Value = 2.0908
$ws.Range("F33").Value = -4.7193
$ws.Range("F34").Value = 22.8807
$ws.Range("F35").Value = 5.3359
$ws.Range("F36").Value = -5.1995
$ws.Range("F37").Value = -5.6238
$ws.Range("F38").Value = -22.595
$ws.Range("F39").Value = 10.8405
$ws.Range("F40").Value = -7.5963
$ws.Range("F41").Value = -4.552
$ws.Range("F42").Value = 22.3098
$ws.Range("F43").Value = 14.0694
$ws.Range("F44").Value = -9.6066
$ws.Range("F45").Value = 27.639
$ws.Range("F46").Value = -6.3484
$ws.Range("F47").Value = -40.5302
$ws.Range("F48").Value = -29.7988
$ws.Range("F49").Value = -24.0791
$ws.Range("F50").Value = -49.1803
$ws.Range("F51").Value = -51.6023
$ws.Range("F52").Value = -34.4756
$ws.Range("F53").Value = -11.5478
$ws.Range("F54").Value = -2.3796
$ws.Range("F55").Value = -15.4382
$ws.Range("F56").Value = -27.6987
$ws.Range("F57").Value = -27.1559
$ws.Range("F58").Value = -2.1585
$ws.Range("F59").Value = -23.0964
$ws.Range("F60").Value = -13.3217
$ws.Range("F61").Value = -8.1496
$ws.Range("F62").Value = -16.0695
$ws.Range("F63").Value = -12.5465
$ws.Range("F64").Value = 47.7264
$ws.Range("F65").Value = -42.4232
$ws.Range("F66").Value = 11.3291
$ws.Range("F67").Value = 14.3746
$ws.Range("F68").Value = 32.6702
$ws.Range("F69").Value = -17.0097
$ws.Range("F70").Value = -13.5162
$ws.Range("F71").Value = 11.4259
$ws.Range("F72").Value = 2.6754
$ws.Range("F73").Value = -11.1574
$ws.Range("F74").Value = -13.2502
$ws.Range("F75").Value = 24.7078
$ws.Range("F76").Value = 53.3554

Write-Host "Update complete"
